$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.99"
$ws.Range("E2").Value = "'0.97%"
$ws.Range("D3").Value = "'35.85"
$ws.Range("E3").Value = "'1.81%"
$ws.Range("D4").Value = "'5.113"
$ws.Range("E4").Value = "'1.38%"
$ws.Range("D5").Value = "'0.08078"
$ws.Range("E5").Value = "'0.99%"
$ws.Range("D6").Value = "'1.923"
$ws.Range("E6").Value = "'-0.63%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'7.757"
$ws.Range("E7").Value = "'0.21%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9271"
$ws.Range("E8").Value = "'0.57%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1350"
$ws.Range("E9").Value = "'4.50%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1906"
$ws.Range("E10").Value = "'3.14%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09179"
$ws.Range("E11").Value = "'-4.90%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03414"
$ws.Range("E12").Value = "'-5.73%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09827"
$ws.Range("E13").Value = "'-0.30%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001441"
$ws.Range("E14").Value = "'3.98%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005756"
$ws.Range("E15").Value = "'-1.98%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.557"
$ws.Range("E16").Value = "'1.49%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.189"
$ws.Range("E17").Value = "'3.61%"
$ws.Range("E18").Value = "'0.61%"
$ws.Range("D19").Value = "'0.3455"
$ws.Range("E19").Value = "'0.70%"
$ws.Range("D20").Value = "'0.1333"
$ws.Range("E20").Value = "'1.73%"
$ws.Range("D21").Value = "'4.897"
$ws.Range("E21").Value = "'-3.00%"
$ws.Range("D22").Value = "'0.2604"
$ws.Range("E22").Value = "'5.59%"
$ws.Range("D23").Value = "'0.04391"
$ws.Range("E23").Value = "'-2.94%"
$ws.Range("E24").Value = "'0.48%"
$ws.Range("D25").Value = "'0.004801"
$ws.Range("E25").Value = "'-0.44%"
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("E26").Value = "'4.07%"
$ws.Range("D27").Value = "'0.0003136"
$ws.Range("E27").Value = "'4.39%"
$ws.Range("D39").Value = "'0.02002"
$ws.Range("E39").Value = "'5.08%"
$ws.Range("D40").Value = "'0.04903"
$ws.Range("E40").Value = "'4.22%"
$ws.Range("D41").Value = "'0.007595"
$ws.Range("E41").Value = "'1.03%"
$ws.Range("D42").Value = "'0.01035"
$ws.Range("E42").Value = "'7.42%"
$ws.Range("D43").Value = "'0.1373"
$ws.Range("E43").Value = "'3.40%"
$ws.Range("D44").Value = "'0.002103"
$ws.Range("E44").Value = "'-0.41%"
$ws.Range("D45").Value = "'0.01095"
$ws.Range("E45").Value = "'0.98%"
$ws.Range("D46").Value = "'0.00006425"
$ws.Range("E46").Value = "'3.06%"
$ws.Range("E47").Value = "'0.07%"
$ws.Range("D49").Value = "'0.001193"
$ws.Range("E49").Value = "'-19.88%"
$ws.Range("E50").Value = "'0.07%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.07%"
